$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1): add new column labels for bank / deposit_type / currency,
# and shift remaining header cells from E onward into the new standard metadata layout.
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Column B (bank) for data rows 2-7
$ws.Range("B2").Value = "臺灣銀行群賢分行"
$ws.Range("B3").Value = "臺灣銀行群賢分行"
$ws.Range("B4").Value = "台北富邦商業銀行襄陽分行"
$ws.Range("B5").Value = "玉山商業銀行雙和分行"
$ws.Range("B6").Value = "臺灣土地銀行雙和分行"
$ws.Range("B7").Value = "美商摩根大通銀行"

# --- Column C (deposit_type) for data rows 2-7
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("C3").Value = "活期儲蓄存款"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("C5").Value = "活期儲蓄存款"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("C7").Value = "活期存款"

# --- Column D (currency) for data rows 2-7
$ws.Range("D2").Value = "新臺幣"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("D7").Value = "美金"

# --- Column E (owner) stays "蕭美琴" for all rows (already there, re-affirm)
$ws.Range("E2").Value = "蕭美琴"
$ws.Range("E3").Value = "蕭美琴"
$ws.Range("E4").Value = "蕭美琴"
$ws.Range("E5").Value = "蕭美琴"
$ws.Range("E6").Value = "蕭美琴"
$ws.Range("E7").Value = "蕭美琴"

# --- Column F (total amount) - move the TWD value that used to sit in column G into F,
# dropping the old secondary USD value entirely.
$ws.Range("F2").Value = 2208108
$ws.Range("F3").Value = 1927695
$ws.Range("F4").Value = 315666
$ws.Range("F5").Value = 667161
$ws.Range("F6").Value = 107566
$ws.Range("F7").Value = 189920

# --- Column G (property_category)
$ws.Range("G2").Value = "deposit"
$ws.Range("G3").Value = "deposit"
$ws.Range("G4").Value = "deposit"
$ws.Range("G5").Value = "deposit"
$ws.Range("G6").Value = "deposit"
$ws.Range("G7").Value = "deposit"

# --- Column H (category)
$ws.Range("H2").Value = "normal"
$ws.Range("H3").Value = "normal"
$ws.Range("H4").Value = "normal"
$ws.Range("H5").Value = "normal"
$ws.Range("H6").Value = "normal"
$ws.Range("H7").Value = "normal"

# --- Column I (date)
$ws.Range("I2").Value = "2012-04-30"
$ws.Range("I3").Value = "2012-04-30"
$ws.Range("I4").Value = "2012-04-30"
$ws.Range("I5").Value = "2012-04-30"
$ws.Range("I6").Value = "2012-04-30"
$ws.Range("I7").Value = "2012-04-30"

# --- Column J (legislator_name)
$ws.Range("J2").Value = "蕭美琴"
$ws.Range("J3").Value = "蕭美琴"
$ws.Range("J4").Value = "蕭美琴"
$ws.Range("J5").Value = "蕭美琴"
$ws.Range("J6").Value = "蕭美琴"
$ws.Range("J7").Value = "蕭美琴"

# --- Column K (legislator_id)
$ws.Range("K2").Value = 981
$ws.Range("K3").Value = 981
$ws.Range("K4").Value = 981
$ws.Range("K5").Value = 981
$ws.Range("K6").Value = 981
$ws.Range("K7").Value = 981

# --- Column L (source_file)
$ws.Range("L2").Value = "tmpcd9a1"
$ws.Range("L3").Value = "tmpcd9a1"
$ws.Range("L4").Value = "tmpcd9a1"
$ws.Range("L5").Value = "tmpcd9a1"
$ws.Range("L6").Value = "tmpcd9a1"
$ws.Range("L7").Value = "tmpcd9a1"

# --- Column M (index)
$ws.Range("M2").Value = 49
$ws.Range("M3").Value = 50
$ws.Range("M4").Value = 51
$ws.Range("M5").Value = 52
$ws.Range("M6").Value = 53
$ws.Range("M7").Value = 54
